$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model names for rows 2..26 (A column), reflecting the re-ordering/renaming
# of the ensemble members after retraining with the new lm.
$names = @(
    "model_15_9_0",
    "model_15_9_22",
    "model_15_9_21",
    "model_15_9_20",
    "model_15_9_19",
    "model_15_9_18",
    "model_15_9_17",
    "model_15_9_16",
    "model_15_9_15",
    "model_15_9_14",
    "model_15_9_13",
    "model_15_9_23",
    "model_15_9_12",
    "model_15_9_10",
    "model_15_9_9",
    "model_15_9_8",
    "model_15_9_7",
    "model_15_9_6",
    "model_15_9_5",
    "model_15_9_4",
    "model_15_9_3",
    "model_15_9_2",
    "model_15_9_1",
    "model_15_9_11",
    "model_15_9_24"
)

# New metric values (columns B..Q), identical for every row after the retrain.
$values = @(
    0.9999384042336427,
    0.9989769827010289,
    0.999968874408921,
    0.9997505669331505,
    0.9999232107551697,
    0.00005749694415287543,
    0.0009549417433199923,
    0.0000598196697030061,
    0.0001492298129745876,
    0.0001045247413387968,
    0.0005240253855816831,
    0.007582673944781975,
    1.000164255376953,
    0.007905484072910305,
    85.52755751348121,
    125.7504597341318
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

Write-Host "Updated rows 2..26 with new model names and metric values"
